$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1215.7428
$ws.Range("I15").Value = 1215.7428
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 3647.2284
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -3478.2284

$ws.Range("H43").Value = 9833.333000000001
$ws.Range("I43").Value = 12500
$ws.Range("J43").Value = 4500
$ws.Range("K43").Value = 12500
$ws.Range("L43").Value = 4500
$ws.Range("M43").Value = -12431
$ws.Range("N43").Value = -4638

$ws.Range("H92").Value = 111464.555
$ws.Range("I92").Value = 111464.555
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 111464.555
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -110216.555

$ws.Range("H138").Value = 4647.281
$ws.Range("I138").Value = 3464.4194
$ws.Range("J138").Value = 6057.615
$ws.Range("K138").Value = 10393.2582
$ws.Range("L138").Value = 18172.845
$ws.Range("M138").Value = -5253.2582
$ws.Range("N138").Value = -28452.845

$ws.Range("H141").Value = 1819.8966
$ws.Range("I141").Value = 1775.7037
$ws.Range("J141").Value = 2416.5
$ws.Range("K141").Value = 5327.1111
$ws.Range("L141").Value = 7249.5
$ws.Range("M141").Value = -147.1111000000001
$ws.Range("N141").Value = -17609.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4351.8774
$ws.Range("I32").Value = 2687.5227
$ws.Range("J32").Value = 18998.2
$ws.Range("K32").Value = 2687.5227
$ws.Range("L32").Value = 18998.2
$ws.Range("M32").Value = -2400.5227
$ws.Range("N32").Value = -19572.2

$ws.Range("H61").Value = 2130
$ws.Range("I61").Value = 2062.5
$ws.Range("J61").Value = 2400
$ws.Range("K61").Value = 2062.5
$ws.Range("L61").Value = 2400
$ws.Range("M61").Value = -1850.5
$ws.Range("N61").Value = -2824

$ws.Range("H74").Value = 1298.375
$ws.Range("I74").Value = 1118.2667
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 1118.2667
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -244.2666999999999
$ws.Range("N74").Value = -5748

$ws.Range("H77").Value = 1298.375
$ws.Range("I77").Value = 1118.2667
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 5591.3335
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -1223.3335
$ws.Range("N77").Value = -28736

$ws.Range("H110").Value = 10011
$ws.Range("I110").Value = 10011
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 10011
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -7966

$ws.Range("H122").Value = 5999.9414
$ws.Range("I122").Value = 5136.273
$ws.Range("J122").Value = 7583.3335
$ws.Range("K122").Value = 15408.819
$ws.Range("L122").Value = 22750.0005
$ws.Range("M122").Value = -12958.819
$ws.Range("N122").Value = -27650.0005

$ws.Range("H132").Value = 1166.5714
$ws.Range("I132").Value = 1166.5714
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3499.7142
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -969.7142000000003
$ws.Range("N132").Value = $null

$ws.Range("H136").Value = 2130
$ws.Range("I136").Value = 2062.5
$ws.Range("J136").Value = 2400
$ws.Range("K136").Value = 6187.5
$ws.Range("L136").Value = 7200
$ws.Range("M136").Value = -3637.5
$ws.Range("N136").Value = -12300

$ws.Range("H138").Value = 74992
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 74992
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 74992
$ws.Range("N138").Value = -85272

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 100
$ws.Range("I8").Value = 100
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 100
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 40

$ws.Range("H94").Value = 1478.3636
$ws.Range("I94").Value = 1502.7
$ws.Range("J94").Value = 1235
$ws.Range("K94").Value = 1502.7
$ws.Range("L94").Value = 1235
$ws.Range("M94").Value = -1051.7
$ws.Range("N94").Value = -2137

$ws.Range("H99").Value = 1867.2
$ws.Range("I99").Value = 1352
$ws.Range("J99").Value = 2640
$ws.Range("K99").Value = 1352
$ws.Range("L99").Value = 2640
$ws.Range("M99").Value = 146
$ws.Range("N99").Value = -5636

$ws.Range("H107").Value = 575.6667
$ws.Range("I107").Value = 622.8
$ws.Range("J107").Value = 340
$ws.Range("K107").Value = 622.8
$ws.Range("L107").Value = 340
$ws.Range("M107").Value = 1297.2
$ws.Range("N107").Value = -4180

$ws.Range("H134").Value = 1793.2142
$ws.Range("I134").Value = 1546.5385
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 4639.6155
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -2104.6155
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 47775.4
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 47775.4
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 47775.4
$ws.Range("N50").Value = -49025.4

$ws.Range("H60").Value = 11116.308
$ws.Range("I60").Value = 11116.308
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 11116.308
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -10605.308
$ws.Range("N60").Value = $null

$ws.Range("H80").Value = 60000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 60000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 60000
$ws.Range("N80").Value = -62246

$ws.Range("H83").Value = 60000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 60000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 180000
$ws.Range("N83").Value = -191232

$ws.Range("H107").Value = 1017.8889
$ws.Range("I107").Value = 715
$ws.Range("J107").Value = 1169.3334
$ws.Range("K107").Value = 715
$ws.Range("L107").Value = 1169.3334
$ws.Range("M107").Value = 1205
$ws.Range("N107").Value = -5009.3334

$ws.Range("H134").Value = 3870
$ws.Range("I134").Value = 3484
$ws.Range("J134").Value = 3998.6667
$ws.Range("K134").Value = 10452
$ws.Range("L134").Value = 11996.0001
$ws.Range("M134").Value = -7917
$ws.Range("N134").Value = -17066.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 100
$ws.Range("I13").Value = 100
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 300
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -132

$ws.Range("H109").Value = 1374.75
$ws.Range("I109").Value = 1500
$ws.Range("J109").Value = 999
$ws.Range("K109").Value = 4500
$ws.Range("L109").Value = 2997
$ws.Range("M109").Value = -3460
$ws.Range("N109").Value = -5077

$ws.Range("H115").Value = 713.8570999999999
$ws.Range("I115").Value = 675
$ws.Range("J115").Value = 729.4
$ws.Range("K115").Value = 2025
$ws.Range("L115").Value = 2188.2
$ws.Range("M115").Value = -850
$ws.Range("N115").Value = -4538.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5786.5713
$ws.Range("I70").Value = 4626.625
$ws.Range("J70").Value = 7333.1665
$ws.Range("K70").Value = 4626.625
$ws.Range("L70").Value = 7333.1665
$ws.Range("M70").Value = -4356.625
$ws.Range("N70").Value = -7873.1665

$ws.Range("H73").Value = 5786.5713
$ws.Range("I73").Value = 4626.625
$ws.Range("J73").Value = 7333.1665
$ws.Range("K73").Value = 4626.625
$ws.Range("L73").Value = 7333.1665
$ws.Range("M73").Value = -3690.625
$ws.Range("N73").Value = -9205.166499999999

$ws.Range("H80").Value = 13083.167
$ws.Range("I80").Value = 6359.8
$ws.Range("J80").Value = 17885.572
$ws.Range("K80").Value = 6359.8
$ws.Range("L80").Value = 17885.572
$ws.Range("M80").Value = -5361.8
$ws.Range("N80").Value = -19881.572

$ws.Range("H83").Value = 13083.167
$ws.Range("I83").Value = 6359.8
$ws.Range("J83").Value = 17885.572
$ws.Range("K83").Value = 31799
$ws.Range("L83").Value = 89427.86
$ws.Range("M83").Value = -26807
$ws.Range("N83").Value = -99411.86

$ws.Range("H113").Value = 2807.1667
$ws.Range("I113").Value = 2210.75
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2210.75
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -40.75
$ws.Range("N113").Value = -8340

$ws.Range("H132").Value = 2831.4614
$ws.Range("I132").Value = 2545.4443
$ws.Range("J132").Value = 3475
$ws.Range("K132").Value = 7636.3329
$ws.Range("L132").Value = 10425
$ws.Range("M132").Value = -5106.3329
$ws.Range("N132").Value = -15485

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5823.1665
$ws.Range("I16").Value = 5823.1665
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 5823.1665
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -5653.1665
$ws.Range("N16").Value = $null

$ws.Range("H46").Value = 3764.3044
$ws.Range("I46").Value = 2921.4285
$ws.Range("J46").Value = 4133.0625
$ws.Range("K46").Value = 2921.4285
$ws.Range("L46").Value = 4133.0625
$ws.Range("M46").Value = -2733.4285
$ws.Range("N46").Value = -4509.0625

$ws.Range("H93").Value = 2266.4443
$ws.Range("I93").Value = 2216.3333
$ws.Range("J93").Value = 2366.6667
$ws.Range("K93").Value = 2216.3333
$ws.Range("L93").Value = 2366.6667
$ws.Range("M93").Value = -968.3332999999998
$ws.Range("N93").Value = -4862.6667

$ws.Range("H100").Value = 2298.3333
$ws.Range("I100").Value = 1947.5
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1947.5
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1406.5
$ws.Range("N100").Value = -4082

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 14999
$ws.Range("I11").Value = 14998.5
$ws.Range("J11").Value = 15000
$ws.Range("K11").Value = 14998.5
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = -14856.5
$ws.Range("N11").Value = -15284

$ws.Range("H14").Value = 10100.25
$ws.Range("I14").Value = 1703
$ws.Range("J14").Value = 18497.5
$ws.Range("K14").Value = 1703
$ws.Range("L14").Value = 18497.5
$ws.Range("M14").Value = -1535
$ws.Range("N14").Value = -18833.5

$ws.Range("H29").Value = 31982.666
$ws.Range("I29").Value = 31982.666
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 31982.666
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -31692.666

$ws.Range("H52").Value = 9259
$ws.Range("I52").Value = 9259
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 9259
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -9033

$ws.Range("H104").Value = 9331.666999999999
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 9331.666999999999
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 9331.666999999999
$ws.Range("N104").Value = -16319.667

$ws.Range("H132").Value = 188807.4
$ws.Range("I132").Value = 188807.4
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 566422.2
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -563892.2

$ws.Range("H136").Value = 1645.4375
$ws.Range("I136").Value = 1645.4375
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4936.3125
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2386.3125
$ws.Range("N136").Value = $null
